# Updated attendance for Willy workshop
$wb = $excel.ActiveWorkbook

# The "Willy Lee" workshop (Oct. 9th) lives on the "2024 - Fall" sheet, row 7.
# Fill in the previously-blank In-Person / Zoom attendance counts.
$ws = $wb.Worksheets.Item("2024 - Fall")
$ws.Activate()

$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 10

# Move the active selection to where the user clicked next (G8), matching
# the saved sheet view state.
$ws.Range("G8").Select()

$excel.Calculate()
